$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) contain values that look numeric (e.g. "91.10") but must
# stay as literal text (matching the source data). Force text format before
# assigning so Excel does not silently coerce them into real numbers and
# strip formatting such as trailing zeros.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.099.54'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.876.12'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.52'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5096'
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3849'
$ws.Range('E8').Value = '  -2.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09112'
$ws.Range('E9').Value = '  -2.36%  '
$ws.Range('E10').Value = '  -1.66%  '
$ws.Range('B11').Value = 'Polkadot'
$ws.Range('C11').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.347'
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.76'
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.874.18'
$ws.Range('E13').Value = '  -1.54%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.211'
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('B15').Value = 'BinanceUSD'
$ws.Range('C15').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.002'
$ws.Range('E15').Value = '  +0.33%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001115'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.10'
$ws.Range('E17').Value = '  -1.69%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06597'
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.17'
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.110'
$ws.Range('E21').Value = '  -1.72%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.116.26'
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.42'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.281'
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('B25').Value = 'LEO'
$ws.Range('C25').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.382'
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.090.20'
$ws.Range('E26').Value = '  -1.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.540'
$ws.Range('E27').Value = '  -2.69%  '
$ws.Range('E28').Value = '  -1.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '157.61'
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.79'
$ws.Range('E30').Value = '  -0.69%  '
$ws.Range('E31').Value = '  -2.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1052'
$ws.Range('E32').Value = '  -1.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.610'
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.596'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.680'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02436'
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06571'
$ws.Range('E37').Value = '  -1.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2177'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.209'
$ws.Range('E39').Value = '  -3.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.263'
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6405'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.914'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6015'
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.275'
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('E48').Value = '  +4.51%  '
$ws.Range('E49').Value = '  -1.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '121.27'
$ws.Range('E50').Value = '  -1.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.70'
$ws.Range('E51').Value = '  +1.67%  '
